# Test plan for savings planed.
#
# Fills in the "Developer" name and the previously-empty Preconditions /
# Method Inputs / Expected Result columns for the SavingsAccount unit
# test plan (rows 7-12), then nudges the view (zoom/selection) and the
# column widths to match the post-edit, content-driven auto-fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Developer name -------------------------------------------------
$ws.Range("C3").Value = "Om Patel"

# --- Test case 1: __init__ / Attributes are set to parameter values. -
$ws.Range("E7").Value = "None "
$ws.Range("F7").Value = "account_number=2004, client_number=2904, balance=1000.0, date_created=date(2024, 10, 5), minimum_balance=50.0"
$ws.Range("G7").Value = "account_number=2004, client_number=2904, balance=1000.0, date_created=date(2024, 10, 5), minimum_balance=50.0"

# --- Test case 2: __init__ / minimum_balance has invalid type. -------
$ws.Range("E8").Value = "None "
$ws.Range("F8").Value = 'account_number=2004, client_number=2904, balance=1000.0, date_created=date(2024, 10, 5), minimum_balance="invalid"'
$ws.Range("G8").Value = "minimum balance is set to be 50.0"

# --- Test case 3: get_service_charges / balance greater than minimum -
$ws.Range("E9").Value = "account_number=2004, client_number=2904, balance=100.0, date_created=date(2024, 10, 5), minimum_balance=50.0"
$ws.Range("F9").Value = "None"
$ws.Range("G9").Value = "Service charge is 0.50"

# --- Test case 4: get_service_charges / balance equal to minimum -----
$ws.Range("E10").Value = "account_number=2004, client_number=2904, balance=50.0, date_created=date(2024, 10, 5), minimum_balance=50.0"
$ws.Range("F10").Value = "None"
$ws.Range("G10").Value = "Service charge is 0.50"

# --- Test case 5: get_service_charges / balance less than minimum ----
$ws.Range("E11").Value = "account_number=2004, client_number=2904, balance=49.99, date_created=date(2024, 10, 5), minimum_balance=50.0"
$ws.Range("F11").Value = "None"
$ws.Range("G11").Value = "Service charge is 1.00"

# --- Test case 6: __str__ / appropriate value returned ----------------
$ws.Range("E12").Value = "account_number=2004, client_number=2904, balance=1000.0, date_created=date(2024, 10, 5), minimum_balance=50.0"
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "String :- Account Number: 2004 Client Number: 2904 Balance: `$1,000.00 Minimum Balance: `$50.00 Account Type: Savings"

# --- Column widths: re-fit now that E/F/G hold long strings -----------
$ws.Columns.Item(2).ColumnWidth = 10.166666666666666
$ws.Columns.Item(3).ColumnWidth = 17.385416666666664
$ws.Columns.Item(4).ColumnWidth = 31.830729166666664
$ws.Columns.Item(5).ColumnWidth = 51.721354166666664
$ws.Columns.Item(6).ColumnWidth = 52.721354166666664
$ws.Columns.Item(7).ColumnWidth = 52.721354166666664

# --- View: zoom out a bit and move the selection onto the new data ----
$excel.ActiveWindow.Zoom = 61
$ws.Range("D7").Select()
